$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Rows 1-3: change the single value in each cell to "0M"
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# 2) Insert 10 new single-value rows right after row 3 (i.e. before the
#    row that currently holds "0" at index 4). Rows.Add(BeforeRow) puts
#    the new row immediately above BeforeRow, so insert the values in
#    reverse order against a fixed anchor to end up in forward order.
$newValues = @("34", "0.00003", "0.00005", "0.00004", "0.00000", "0.00003", "0.00004", "0.00004", "0.00122", "100.0")
$anchorRow = $t.Rows.Item(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($anchorRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
}

# 3) Collapse the last three tab-separated rows down to a single value
#    each, matching what the first three rows used to contain.
$lastIndex = $t.Rows.Count
$t.Rows.Item($lastIndex - 2).Cells.Item(1).Range.Text = "99.99"
$t.Rows.Item($lastIndex - 1).Cells.Item(1).Range.Text = "0"
$t.Rows.Item($lastIndex).Cells.Item(1).Range.Text = "18"
